# Adapt column header formatting to respective input file names (#7)
#
# 1. Rename the "_old" / "_new" header-name suffixes used in the diff sheet
#    to the respective format-version suffixes ("_FV2310" / "_FV2404").
# 2. Wrap the data range in an Excel Table ("Table1").
# 3. Freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) -----------------------------------------
# A1:J1 used the "_old" suffix -> "_FV2310"
# K1 ("diff") is unchanged
# L1:U1 used the "_new" suffix -> "_FV2404"

$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($baseNames[$i])_FV2310"
}

$ws.Cells.Item(1, 11).Value = "diff"

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = "$($baseNames[$i])_FV2404"
}

# --- 2. Freeze the header row -----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range into an Excel Table -----------------------------
$tableRange = $ws.Range("A1:U58")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""
